{"js": "// The document ends with a \"Requisitos\" section whose last requirement\n// line (\"LOT2039: ...\") was followed by an empty paragraph, a page-break\n// paragraph, and a site-footer \"\u00a9 2020 ...\" paragraph. The edit removes\n// those three paragraphs, leaving the LOT2039 line directly followed by\n// the (unrelated) trailing empty / page-break paragraphs that close the\n// document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph by its literal text (robust to any\n// surrounding paragraphs shifting around).\nconst anchorText =\n  \"LOT2039: Estrutura e Qu\u00edmica de Materiais Lignocelul\u00f3sicos (Requisito)\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the LOT2039 requisito paragraph\");\n}\n\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// The three paragraphs to delete are the ones right after the anchor:\n// an empty paragraph, an empty page-break paragraph, and the copyright\n// paragraph itself. Confirm shape before deleting, then delete from the\n// end backwards so indices of earlier paragraphs stay valid.\nconst toDelete = [];\nif (\n  anchorIndex + 3 <= items.length - 1 &&\n  items[anchorIndex + 1].text === \"\" &&\n  items[anchorIndex + 2].text === \"\" &&\n  items[anchorIndex + 3].text === copyrightText\n) {\n  toDelete.push(items[anchorIndex + 1], items[anchorIndex + 2], items[anchorIndex + 3]);\n} else {\n  // Fallback: delete the copyright paragraph and the up-to-two empty\n  // paragraphs that immediately precede it (between it and the anchor).\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === copyrightText) {\n      let start = i;\n      while (start - 1 > anchorIndex && items[start - 1].text === \"\") {\n        start--;\n      }\n      for (let j = start; j <= i; j++) toDelete.push(items[j]);\n      break;\n    }\n  }\n}\n\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The document ends with a \"Requisitos\" section whose last requirement\n# line (\"LOT2039: ...\") was followed by an empty paragraph, a page-break\n# paragraph, and a site-footer \"\u00a9 2020 ...\" paragraph. The edit removes\n# those three paragraphs, leaving the LOT2039 line directly followed by\n# the (unrelated) trailing empty / page-break paragraphs that close the\n# document.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOT2039: Estrutura e Qu\u00edmica de Materiais Lignocelul\u00f3sicos (Requisito)\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\nfunction Clean-ParaText($t) {\n    return $t.TrimEnd([char]13, [char]7)\n}\n\n# Locate the anchor paragraph by its literal text (robust to any\n# surrounding paragraphs shifting around).\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ((Clean-ParaText $d.Paragraphs.Item($i).Range.Text) -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the LOT2039 requisito paragraph\"\n}\n\n$deleteFrom = -1\n$deleteTo = -1\n\nif (($anchorIndex + 3) -le $d.Paragraphs.Count `\n    -and (Clean-ParaText $d.Paragraphs.Item($anchorIndex + 1).Range.Text) -eq \"\" `\n    -and (Clean-ParaText $d.Paragraphs.Item($anchorIndex + 2).Range.Text) -eq \"\" `\n    -and (Clean-ParaText $d.Paragraphs.Item($anchorIndex + 3).Range.Text) -eq $copyrightText) {\n    $deleteFrom = $anchorIndex + 1\n    $deleteTo = $anchorIndex + 3\n}\nelse {\n    # Fallback: locate the copyright paragraph and include the (up to two)\n    # empty paragraphs immediately preceding it, down to the anchor.\n    for ($i = $anchorIndex + 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ((Clean-ParaText $d.Paragraphs.Item($i).Range.Text) -eq $copyrightText) {\n            $start = $i\n            while (($start - 1) -gt $anchorIndex -and (Clean-ParaText $d.Paragraphs.Item($start - 1).Range.Text) -eq \"\") {\n                $start--\n            }\n            $deleteFrom = $start\n            $deleteTo = $i\n            break\n        }\n    }\n}\n\nif ($deleteFrom -eq -1) {\n    throw \"Could not locate the paragraphs to delete\"\n}\n\n$rangeStart = $d.Paragraphs.Item($deleteFrom).Range.Start\n$rangeEnd = $d.Paragraphs.Item($deleteTo).Range.End\n$r = $d.Range($rangeStart, $rangeEnd)\n$r.Delete()\n"}
